$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Labels for the new matrix block (row 85) ---
# Set "Matrix " first so it lands at shared-string index 21, then
# "Matrix Transpose" at index 22 (matches target sharedStrings order).
$ws.Range("K85").Value = "Matrix "
$ws.Range("E85").Value = "Matrix Transpose"

# --- Highlight the 2x5 "Matrix" results block (A86:B90) in yellow ---
$ws.Range("A86:B90").Interior.Color = 65535

# --- Row 86: identity-ish helper row + transpose seeds ---
$ws.Range("E86").Value = 1
$ws.Range("F86").Value = 1
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 1
$ws.Range("I86").Value = 1
$ws.Range("K86").Value = 1
$ws.Range("L86").Formula = "=A86"
$ws.Range("M86").Formula = "=B86"

# --- Row 87 ---
$ws.Range("E87").Formula = "=`$L`$86"
$ws.Range("F87").Formula = "=`$L`$87"
$ws.Range("G87").Formula = "=`$L`$88"
$ws.Range("H87").Formula = "=`$L`$89"
$ws.Range("I87").Formula = "=`$L`$90"
$ws.Range("K87").Value = 1
$ws.Range("L87").Formula = "=A87"
$ws.Range("M87").Formula = "=B87"

# --- Row 88 ---
$ws.Range("E88").Formula = "=`$L`$86"
$ws.Range("F88").Formula = "=`$L`$87"
$ws.Range("G88").Formula = "=`$L`$88"
$ws.Range("H88").Formula = "=`$L`$89"
$ws.Range("I88").Formula = "=`$L`$90"
$ws.Range("K88").Value = 1
$ws.Range("L88").Formula = "=A88"
$ws.Range("M88").Formula = "=B88"

# --- Row 89 ---
$ws.Range("K89").Value = 1
$ws.Range("L89").Formula = "=A89"
$ws.Range("M89").Formula = "=B89"

# --- Row 90 ---
$ws.Range("K90").Value = 1
$ws.Range("L90").Formula = "=A90"
$ws.Range("M90").Formula = "=B90"

# --- View state: selection moved to the new block ---
$ws.Range("K86").Select() | Out-Null
